$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 11).Value2 = 15
$ws.Cells.Item(2, 14).Value2 = 1.62
$ws.Cells.Item(2, 15).Value2 = 2.25

# Row 6
$ws.Cells.Item(6, 12).Value2 = 1.3
$ws.Cells.Item(6, 13).Value2 = 3.5
$ws.Cells.Item(6, 14).Value2 = 2.04
$ws.Cells.Item(6, 15).Value2 = 1.86

# Row 8
$ws.Cells.Item(8, 8).Value2 = 5.25
$ws.Cells.Item(8, 9).Value2 = 12
$ws.Cells.Item(8, 14).Value2 = 1.53
$ws.Cells.Item(8, 15).Value2 = 2.4
$ws.Cells.Item(8, 18).Value2 = 2
$ws.Cells.Item(8, 19).Value2 = 1.73
$ws.Cells.Item(8, 27).Value2 = 10
$ws.Cells.Item(8, 28).Value2 = 21

# Row 9
$ws.Cells.Item(9, 7).Value2 = 3
$ws.Cells.Item(9, 8).Value2 = 3.3
$ws.Cells.Item(9, 9).Value2 = 2.32
$ws.Cells.Item(9, 11).Value2 = 7.7
$ws.Cells.Item(9, 16).Value2 = 1.42
$ws.Cells.Item(9, 17).Value2 = 2.7
$ws.Cells.Item(9, 19).Value2 = 2.02
$ws.Cells.Item(9, 20).Value2 = 9.5
$ws.Cells.Item(9, 21).Value2 = 17
$ws.Cells.Item(9, 22).Value2 = 11
$ws.Cells.Item(9, 23).Value2 = 40
$ws.Cells.Item(9, 24).Value2 = 26
$ws.Cells.Item(9, 25).Value2 = 35
$ws.Cells.Item(9, 26).Value2 = 7.7
$ws.Cells.Item(9, 27).Value2 = 6.6
$ws.Cells.Item(9, 31).Value2 = 7.8
$ws.Cells.Item(9, 32).Value2 = 12
$ws.Cells.Item(9, 34).Value2 = 25
$ws.Cells.Item(9, 35).Value2 = 20

# Row 15
$ws.Cells.Item(15, 8).Value2 = 2.8
$ws.Cells.Item(15, 11).Value2 = 5
$ws.Cells.Item(15, 17).Value2 = 2.22
$ws.Cells.Item(15, 18).Value2 = 2.2
$ws.Cells.Item(15, 19).Value2 = 1.6
$ws.Cells.Item(15, 20).Value2 = 7.9
$ws.Cells.Item(15, 22).Value2 = 14
$ws.Cells.Item(15, 26).Value2 = 5
$ws.Cells.Item(15, 27).Value2 = 5.7
$ws.Cells.Item(15, 29).Value2 = 150
$ws.Cells.Item(15, 31).Value2 = 5.3
$ws.Cells.Item(15, 32).Value2 = 8.75
$ws.Cells.Item(15, 33).Value2 = 9.5

# Row 18
$ws.Cells.Item(18, 12).Value2 = 1.3
$ws.Cells.Item(18, 13).Value2 = 3.5
$ws.Cells.Item(18, 14).Value2 = 1.98
$ws.Cells.Item(18, 15).Value2 = 1.88

# Row 20
$ws.Cells.Item(20, 9).Value2 = 2.2
$ws.Cells.Item(20, 19).Value2 = 2.35
$ws.Cells.Item(20, 21).Value2 = 18
$ws.Cells.Item(20, 24).Value2 = 21
$ws.Cells.Item(20, 25).Value2 = 23
$ws.Cells.Item(20, 26).Value2 = 14.5
$ws.Cells.Item(20, 31).Value2 = 10.75
$ws.Cells.Item(20, 32).Value2 = 13
$ws.Cells.Item(20, 34).Value2 = 23

# Row 21
$ws.Cells.Item(21, 10).Value2 = 1.08
$ws.Cells.Item(21, 11).Value2 = 8
$ws.Cells.Item(21, 14).Value2 = 2.1
$ws.Cells.Item(21, 15).Value2 = 1.7
$ws.Cells.Item(21, 21).Value2 = 6
$ws.Cells.Item(21, 22).Value2 = 9
$ws.Cells.Item(21, 26).Value2 = 8
$ws.Cells.Item(21, 31).Value2 = 15

# Row 22
$ws.Cells.Item(22, 7).Value2 = 3.2
$ws.Cells.Item(22, 8).Value2 = 3
$ws.Cells.Item(22, 9).Value2 = 2.38
$ws.Cells.Item(22, 22).Value2 = 12
$ws.Cells.Item(22, 24).Value2 = 26
$ws.Cells.Item(22, 35).Value2 = 19

# Row 24
$ws.Cells.Item(24, 14).Value2 = 1.7
$ws.Cells.Item(24, 15).Value2 = 2.1

# Row 25
$ws.Cells.Item(25, 7).Value2 = 4.2
$ws.Cells.Item(25, 9).Value2 = 1.95
$ws.Cells.Item(25, 10).Value2 = 1.06
$ws.Cells.Item(25, 11).Value2 = 10
$ws.Cells.Item(25, 18).Value2 = 1.75
$ws.Cells.Item(25, 19).Value2 = 2
$ws.Cells.Item(25, 22).Value2 = 15
$ws.Cells.Item(25, 24).Value2 = 34
$ws.Cells.Item(25, 25).Value2 = 41
$ws.Cells.Item(25, 31).Value2 = 7.5
$ws.Cells.Item(25, 32).Value2 = 9.5

# Row 27
$ws.Cells.Item(27, 7).Value2 = 3.8
$ws.Cells.Item(27, 8).Value2 = 3.75
$ws.Cells.Item(27, 9).Value2 = 1.9
$ws.Cells.Item(27, 11).Value2 = 12
$ws.Cells.Item(27, 16).Value2 = 1.36
$ws.Cells.Item(27, 17).Value2 = 3
$ws.Cells.Item(27, 29).Value2 = 51
$ws.Cells.Item(27, 32).Value2 = 9.5
$ws.Cells.Item(27, 34).Value2 = 17

# Row 28
$ws.Cells.Item(28, 12).Value2 = 1.22
$ws.Cells.Item(28, 13).Value2 = 4
$ws.Cells.Item(28, 14).Value2 = 1.8
$ws.Cells.Item(28, 15).Value2 = 2

# Row 29
$ws.Cells.Item(29, 7).Value2 = 1.44
$ws.Cells.Item(29, 8).Value2 = 4.5
$ws.Cells.Item(29, 9).Value2 = 7
$ws.Cells.Item(29, 34).Value2 = 81
$ws.Cells.Item(29, 35).Value2 = 51

# Row 30
$ws.Cells.Item(30, 7).Value2 = 3.4
$ws.Cells.Item(30, 8).Value2 = 3.5
$ws.Cells.Item(30, 9).Value2 = 2.05
$ws.Cells.Item(30, 10).Value2 = 1.06
$ws.Cells.Item(30, 11).Value2 = 10
$ws.Cells.Item(30, 14).Value2 = 2.15
$ws.Cells.Item(30, 15).Value2 = 1.67
$ws.Cells.Item(30, 22).Value2 = 12
$ws.Cells.Item(30, 24).Value2 = 29
$ws.Cells.Item(30, 30).Value2 = 351
$ws.Cells.Item(30, 34).Value2 = 19
$ws.Cells.Item(30, 36).Value2 = 29

# Row 33
$ws.Cells.Item(33, 16).Value2 = 1.57
$ws.Cells.Item(33, 17).Value2 = 2.32
$ws.Cells.Item(33, 28).Value2 = 19

# Row 36
$ws.Cells.Item(36, 14).Value2 = 1.95
$ws.Cells.Item(36, 15).Value2 = 1.85

# Row 37
$ws.Cells.Item(37, 22).Value2 = 9.5

# Row 38
$ws.Cells.Item(38, 7).Value2 = 1.67
$ws.Cells.Item(38, 8).Value2 = 4.33
$ws.Cells.Item(38, 9).Value2 = 4.2
$ws.Cells.Item(38, 14).Value2 = 1.57
$ws.Cells.Item(38, 15).Value2 = 2.35
$ws.Cells.Item(38, 16).Value2 = 1.29
$ws.Cells.Item(38, 17).Value2 = 3.5
$ws.Cells.Item(38, 26).Value2 = 17
$ws.Cells.Item(38, 27).Value2 = 8.5
$ws.Cells.Item(38, 31).Value2 = 15
$ws.Cells.Item(38, 32).Value2 = 23
$ws.Cells.Item(38, 33).Value2 = 13
$ws.Cells.Item(38, 34).Value2 = 41
$ws.Cells.Item(38, 35).Value2 = 29
$ws.Cells.Item(38, 36).Value2 = 29

# Row 39
$ws.Cells.Item(39, 11).Value2 = 9.25
$ws.Cells.Item(39, 12).Value2 = 1.17
$ws.Cells.Item(39, 13).Value2 = 4.45
$ws.Cells.Item(39, 14).Value2 = 1.52
$ws.Cells.Item(39, 15).Value2 = 2.35
$ws.Cells.Item(39, 16).Value2 = 1.28
$ws.Cells.Item(39, 17).Value2 = 3.35
$ws.Cells.Item(39, 18).Value2 = 1.57
$ws.Cells.Item(39, 19).Value2 = 2.27
$ws.Cells.Item(39, 20).Value2 = 9.75
$ws.Cells.Item(39, 21).Value2 = 9.75
$ws.Cells.Item(39, 23).Value2 = 14
$ws.Cells.Item(39, 25).Value2 = 19
$ws.Cells.Item(39, 26).Value2 = 9.25
$ws.Cells.Item(39, 28).Value2 = 13
$ws.Cells.Item(39, 29).Value2 = 45
$ws.Cells.Item(39, 30).Value2 = 250
$ws.Cells.Item(39, 31).Value2 = 16.5
$ws.Cells.Item(39, 32).Value2 = 27
$ws.Cells.Item(39, 33).Value2 = 13.5
$ws.Cells.Item(39, 35).Value2 = 32
$ws.Cells.Item(39, 36).Value2 = 32

# Row 40
$ws.Cells.Item(40, 7).Value2 = 8.75
$ws.Cells.Item(40, 8).Value2 = 5.7
$ws.Cells.Item(40, 10).Value2 = 1.02
$ws.Cells.Item(40, 11).Value2 = 10
$ws.Cells.Item(40, 12).Value2 = 1.13
$ws.Cells.Item(40, 13).Value2 = 5.1
$ws.Cells.Item(40, 14).Value2 = 1.42
$ws.Cells.Item(40, 15).Value2 = 2.65
$ws.Cells.Item(40, 16).Value2 = 1.24
$ws.Cells.Item(40, 17).Value2 = 3.65
$ws.Cells.Item(40, 18).Value2 = 1.82
$ws.Cells.Item(40, 19).Value2 = 1.9
$ws.Cells.Item(40, 20).Value2 = 29
$ws.Cells.Item(40, 24).Value2 = 90
$ws.Cells.Item(40, 25).Value2 = 70
$ws.Cells.Item(40, 26).Value2 = 10
$ws.Cells.Item(40, 27).Value2 = 11.75
$ws.Cells.Item(40, 30).Value2 = 500
$ws.Cells.Item(40, 31).Value2 = 9.75
$ws.Cells.Item(40, 32).Value2 = 7.4
$ws.Cells.Item(40, 33).Value2 = 9
$ws.Cells.Item(40, 35).Value2 = 10
$ws.Cells.Item(40, 36).Value2 = 23

# Row 41
$ws.Cells.Item(41, 7).Value2 = 2.77
$ws.Cells.Item(41, 8).Value2 = 3.45
$ws.Cells.Item(41, 9).Value2 = 2.32
$ws.Cells.Item(41, 11).Value2 = 8.5
$ws.Cells.Item(41, 12).Value2 = 1.2
$ws.Cells.Item(41, 13).Value2 = 4.05
$ws.Cells.Item(41, 14).Value2 = 1.6
$ws.Cells.Item(41, 15).Value2 = 2.18
$ws.Cells.Item(41, 16).Value2 = 1.32
$ws.Cells.Item(41, 17).Value2 = 3.1
$ws.Cells.Item(41, 18).Value2 = 1.5
$ws.Cells.Item(41, 19).Value2 = 2.42
$ws.Cells.Item(41, 26).Value2 = 8.5
$ws.Cells.Item(41, 27).Value2 = 6.9
$ws.Cells.Item(41, 28).Value2 = 11.25
$ws.Cells.Item(41, 29).Value2 = 37
$ws.Cells.Item(41, 30).Value2 = 200
$ws.Cells.Item(41, 31).Value2 = 11.25
$ws.Cells.Item(41, 32).Value2 = 14.5
$ws.Cells.Item(41, 34).Value2 = 26
$ws.Cells.Item(41, 35).Value2 = 16.5
$ws.Cells.Item(41, 36).Value2 = 20
